$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: bonus hour entry date moves from 2022-03-01 to 2022-03-02, and the
# row grows taller to fit the (unchanged) wrapped description text.
$ws.Rows.Item(20).RowHeight = 28.45
$ws.Range("C20").Value = 44622

# Row 21: new bonus-hour log entry describing further MVVM refactoring work.
$ws.Rows.Item(21).RowHeight = 43
$ws.Range("A21").Value = "Just some more refactoring"
$ws.Range("B21").Value = 1
$ws.Range("C21").Value = 44622
$ws.Range("D21").Value = "Added an infoViewModel for consistency. Took out the houses property in HousesViewController and made use of the HouseManager.houses. Same for chosenHouse in DetailViewController. Now I think it’s proper use of the MVVM pattern."
